# add export folder path & promotion price for WAT
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- "MAN product": the leading index column (bold + bordered 0/1 markers)
#     goes away entirely, shifting every other column one to the left - the
#     "MAN ID"/572123/433243 text slides from column B into column A and the
#     date stamps slide from column H into column G. ---
$ws1.Columns.Item(1).Delete()

$ws1.Range("A4:B4").Select() | Out-Null

# --- Add the new "WAT url" sheet right after "MAN product" with the two
#     promo-page links. ---
$wsWat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$wsWat.Name = "WAT url"

$url1 = "https://www.watsons.com.hk/%E8%AD%B7%E8%88%92%E5%AF%B6radiant%E6%97%A5%E7%94%A824cm-9%E7%89%87/p/BP_287456"
$url2 = "https://www.watsons.com.hk/%E6%BB%8B%E6%BD%A4%E8%82%B2%E9%AB%AE%E7%B2%BE%E8%8F%AF%E7%B4%A0/p/BP_266919"

$wsWat.Range("A1").Value2 = "WAT url"
$wsWat.Range("A2").Value2 = $url1
$wsWat.Range("A3").Value2 = $url2

$wsWat.Hyperlinks.Add($wsWat.Range("A2"), $url1) | Out-Null
$wsWat.Hyperlinks.Add($wsWat.Range("A3"), $url2) | Out-Null

$wsWat.Range("A3").Select() | Out-Null
$wsWat.Activate() | Out-Null
